$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

# Insert a new column at C, shifting existing C (dbExcel) and D (WebExcel) to D and E
$ws.Columns.Item(3).Insert()

# Header for new column C
$ws.Range("C1").Value2 = "caseDetailQuery"

# Set the width for the newly inserted column C to match column A's width
$ws.Columns.Item(3).ColumnWidth = 74.92

# Content for new C2 cell (query used for reading the case detail table)
$query = @'
MATCH (f:file)-[*]->(c:case) WITH DISTINCT(f) AS f, c MATCH (f)-->(parent) WHERE c.case_id IN ['NCATS-COP01CCB010072'] RETURN f.file_name AS `File Name` ,f.file_type AS `File Type`,head(labels(parent)) AS `Association`, f.file_description AS `Description`,f.file_format AS Format,((f.file_size)/1024) AS Size
'@
$ws.Range("C2").Value2 = $query

# Apply the wrap-text formatting (same as A2) to B2 and C2
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true

# Update the selection / active cell to B2 as in the target view
[void]$ws.Range("B2").Select()
